$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 489.5
$ws.Range("I28").Value = 482.4
$ws.Range("J28").Value = 525
$ws.Range("K28").Value = 482.4
$ws.Range("L28").Value = 525
$ws.Range("M28").Value = 2.600000000000023
$ws.Range("N28").Value = -1495
$ws.Range("H29").Value = 2161.875
$ws.Range("I29").Value = 1897.5
$ws.Range("J29").Value = 2250
$ws.Range("K29").Value = 5692.5
$ws.Range("L29").Value = 6750
$ws.Range("M29").Value = -5411.5
$ws.Range("N29").Value = -7312
$ws.Range("H43").Value = 2425.7144
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2425.7144
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2425.7144
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2563.7144
$ws.Range("H48").Value = 1973
$ws.Range("I48").Value = 300
$ws.Range("J48").Value = 2809.5
$ws.Range("K48").Value = 900
$ws.Range("L48").Value = 8428.5
$ws.Range("M48").Value = -608
$ws.Range("N48").Value = -9012.5
$ws.Range("H56").Value = 1973
$ws.Range("I56").Value = 300
$ws.Range("J56").Value = 2809.5
$ws.Range("K56").Value = 900
$ws.Range("L56").Value = 8428.5
$ws.Range("M56").Value = -366
$ws.Range("N56").Value = -9496.5
$ws.Range("H98").Value = 2619.1667
$ws.Range("I98").Value = 2503.36
$ws.Range("J98").Value = 3198.2
$ws.Range("K98").Value = 2503.36
$ws.Range("L98").Value = 3198.2
$ws.Range("M98").Value = -1005.36
$ws.Range("N98").Value = -6194.2
$ws.Range("H107").Value = 450.9565
$ws.Range("I107").Value = 370.33334
$ws.Range("K107").Value = 370.33334
$ws.Range("M107").Value = 1549.66666
$ws.Range("H116").Value = 5903.385
$ws.Range("I116").Value = 5974.857
$ws.Range("K116").Value = 5974.857
$ws.Range("M116").Value = -2532.857
$ws.Range("H122").Value = 2619.1667
$ws.Range("I122").Value = 2503.36
$ws.Range("J122").Value = 3198.2
$ws.Range("K122").Value = 7510.08
$ws.Range("L122").Value = 9594.599999999999
$ws.Range("M122").Value = -5060.08
$ws.Range("N122").Value = -14494.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 7191.5713
$ws.Range("I28").Value = 7191.5713
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 7191.5713
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -6999.5713
$ws.Range("N28").ClearContents()
$ws.Range("H74").Value = 1083.421
$ws.Range("I74").Value = 1098.1875
$ws.Range("K74").Value = 1098.1875
$ws.Range("M74").Value = -224.1875
$ws.Range("H77").Value = 1083.421
$ws.Range("I77").Value = 1098.1875
$ws.Range("K77").Value = 5490.9375
$ws.Range("M77").Value = -1122.9375
$ws.Range("H99").Value = 7191.5713
$ws.Range("I99").Value = 7191.5713
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 7191.5713
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -4196.5713
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2007.2
$ws.Range("I20").Value = 1891.5
$ws.Range("J20").Value = 2084.3333
$ws.Range("K20").Value = 1891.5
$ws.Range("L20").Value = 2084.3333
$ws.Range("M20").Value = -1644.5
$ws.Range("N20").Value = -2578.3333
$ws.Range("H99").Value = 2477.1428
$ws.Range("I99").Value = 1320
$ws.Range("J99").Value = 3120
$ws.Range("K99").Value = 1320
$ws.Range("L99").Value = 3120
$ws.Range("M99").Value = 178
$ws.Range("N99").Value = -6116

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 19234932
$ws.Range("I58").Value = 3194.5881
$ws.Range("J58").Value = 55561544
$ws.Range("K58").Value = 3194.5881
$ws.Range("L58").Value = 55561544
$ws.Range("M58").Value = -2991.5881
$ws.Range("N58").Value = -55561950
$ws.Range("H74").Value = 18422.75
$ws.Range("J74").Value = 18422.75
$ws.Range("L74").Value = 18422.75
$ws.Range("N74").Value = -20170.75
$ws.Range("H77").Value = 18422.75
$ws.Range("J77").Value = 18422.75
$ws.Range("L77").Value = 55268.25
$ws.Range("N77").Value = -64004.25
$ws.Range("H105").Value = 2806.6
$ws.Range("J105").Value = 4375
$ws.Range("L105").Value = 4375
$ws.Range("N105").Value = -7869
$ws.Range("H134").Value = 2668.9333
$ws.Range("I134").Value = 953.4
$ws.Range("J134").Value = 6100
$ws.Range("K134").Value = 2860.2
$ws.Range("L134").Value = 18300
$ws.Range("M134").Value = -325.1999999999998
$ws.Range("N134").Value = -23370
$ws.Range("H136").Value = 19234932
$ws.Range("I136").Value = 3194.5881
$ws.Range("J136").Value = 55561544
$ws.Range("K136").Value = 9583.764299999999
$ws.Range("L136").Value = 166684632
$ws.Range("M136").Value = -7033.764299999999
$ws.Range("N136").Value = -166689732

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 9996.857
$ws.Range("I87").Value = 1992.6666
$ws.Range("K87").Value = 5977.9998
$ws.Range("M87").Value = -4729.9998
$ws.Range("H90").Value = 9996.857
$ws.Range("I90").Value = 1992.6666
$ws.Range("K90").Value = 17933.9994
$ws.Range("M90").Value = -11693.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3211.738
$ws.Range("I132").Value = 2829.8333
$ws.Range("J132").Value = 4166.5
$ws.Range("K132").Value = 8489.499899999999
$ws.Range("L132").Value = 12499.5
$ws.Range("M132").Value = -5959.499899999999
$ws.Range("N132").Value = -17559.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1578.9166
$ws.Range("I93").Value = 771.44446
$ws.Range("K93").Value = 771.44446
$ws.Range("M93").Value = 476.55554

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 28838.066
$ws.Range("I82").Value = 12273
$ws.Range("J82").Value = 30021.285
$ws.Range("K82").Value = 12273
$ws.Range("L82").Value = 30021.285
$ws.Range("M82").Value = -11890
$ws.Range("N82").Value = -30787.285
$ws.Range("H85").Value = 28838.066
$ws.Range("I85").Value = 12273
$ws.Range("J85").Value = 30021.285
$ws.Range("K85").Value = 12273
$ws.Range("L85").Value = 30021.285
$ws.Range("M85").Value = -10947
$ws.Range("N85").Value = -32673.285
$ws.Range("H107").Value = 704.88
$ws.Range("I107").Value = 297.1
$ws.Range("J107").Value = 2336
$ws.Range("K107").Value = 891.3000000000001
$ws.Range("L107").Value = 7008
$ws.Range("M107").Value = 1028.7
$ws.Range("N107").Value = -10848
$ws.Range("H132").Value = 141552.9
$ws.Range("I132").Value = 170298.4
$ws.Range("J132").Value = 11092.538
$ws.Range("K132").Value = 510895.2
$ws.Range("L132").Value = 33277.614
$ws.Range("M132").Value = -508365.2
$ws.Range("N132").Value = -38337.614
$ws.Range("H136").Value = 1226.1818
$ws.Range("I136").Value = 577.6539
$ws.Range("J136").Value = 3635
$ws.Range("K136").Value = 1732.9617
$ws.Range("L136").Value = 10905
$ws.Range("M136").Value = 817.0382999999999
$ws.Range("N136").Value = -16005

Write-Host "Applied scheduled runner updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
